$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") serial date value moves from 45183 (2023-09-14)
# to 45184 (2023-09-15) for every data row (rows 2-17).
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 3).Value = 45184
}
